# Appends a new "2025-03-21" row (row 20) to each price sheet in the
# workbook, carrying forward the same price that was recorded on the
# previous day (row 19), matching the author's "Updated Argent prices"
# commit.

$wb = $excel.ActiveWorkbook

$sheetPrices = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "43"
    "N-type Wafer"               = "1.19"
    "Cell Topcon 183mm"          = "0.298"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,444"
    "Silver Busbar front-side"   = "8,150"
    "Silver finger front-side"   = "8,200"
    "USD_CNY"                    = "7.2481"
}

foreach ($sheetName in $sheetPrices.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $price = $sheetPrices[$sheetName]

    $dateCell = $ws.Cells.Item(20, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2025-03-21"
    $dateCell.Style = "Normal"

    $priceCell = $ws.Cells.Item(20, 2)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price
    $priceCell.Style = "Normal"
}
